# Revenue and Expenses Completed.
# - Two members' Role changed from SUPER_ADMIN to USER.
# - The "RefferenceNumber" column had several cells that accidentally held a
#   full MemberModel(...) debug string; the stray "referenceMember=null, "
#   fragment is stripped from all of them, and the one describing member
#   no=7 (Manesh Micheal) also has its embedded role updated to USER to
#   match that member's corrected role.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Role column (O) corrections -------------------------------------------------
# Row 3  -> member no=2  (ജെയ്‌സൺ  കൊളുത്തുവള്ളി)
# Row 8  -> member no=7  (മനേഷ് മൈക്കിൽ)
$ws.Range("O3").Value = "USER"
$ws.Range("O8").Value = "USER"

# --- RefferenceNumber column (G) stray debug-string cleanup ----------------------
$ws.Range("G14").Value = "MemberModel(no=6, memberName=ലൂയിസ് ജോസഫ്, emailId=luisjoseph6551@gmail.com, password=bHVpc2pvc2VwaDY1NTE=, memberType=PRIMARY, role=USER, memberDOB=2025-05-17, mobileNumber=9947886919, address1=null, place=null, taluk=Cherthala, district=Alappuzha, state=Kerala, pincode=688537)"
$ws.Range("G15").Value = "MemberModel(no=6, memberName=ലൂയിസ് ജോസഫ്, emailId=luisjoseph6551@gmail.com, password=bHVpc2pvc2VwaDY1NTE=, memberType=PRIMARY, role=USER, memberDOB=2025-05-17, mobileNumber=9947886919, address1=null, place=null, taluk=Cherthala, district=Alappuzha, state=Kerala, pincode=688537)"
$ws.Range("G16").Value = "MemberModel(no=6, memberName=ലൂയിസ് ജോസഫ്, emailId=luisjoseph6551@gmail.com, password=bHVpc2pvc2VwaDY1NTE=, memberType=PRIMARY, role=USER, memberDOB=2025-05-17, mobileNumber=9947886919, address1=null, place=null, taluk=Cherthala, district=Alappuzha, state=Kerala, pincode=688537)"
$ws.Range("G17").Value = "MemberModel(no=5, memberName=ജിജോ ജോൺ, emailId=jijojohn777@gmail.com, password=amlqb2pvaG43Nzc=, memberType=PRIMARY, role=SUPER_ADMIN, memberDOB=2025-12-20, mobileNumber=9895415828, address1=null, place=null, taluk=Cherthala, district=Alappuzha, state=Kerala, pincode=688537)"
$ws.Range("G18").Value = "MemberModel(no=12, memberName=സുജിത്ത് കൃഷ്ണാ, emailId=sujithskrishna@gmail.com, password=c3VqaXRoc2tyaXNobmE=, memberType=PRIMARY, role=SUPER_ADMIN, memberDOB=2025-07-06, mobileNumber=9845375129, address1=Kaithavalappil House, place=Parayakad PO, taluk=Cherthala, district=Alappuzha, state=Kerala, pincode=688540)"
$ws.Range("G19").Value = "MemberModel(no=3, memberName=ജിജിൻ ഫ്രാൻസിസ്, emailId=jijinfrancisanna@gmail.com, password=amlqaW5mcmFuY2lzYW5uYQ==, memberType=PRIMARY, role=SUPER_ADMIN, memberDOB=2025-11-18, mobileNumber=9846417927, address1=null, place=null, taluk=Cherthala, district=Alappuzha, state=Kerala, pincode=688537)"
$ws.Range("G20").Value = "MemberModel(no=10, memberName=സരിത റഫയൽ, emailId=saritharaphel@gmail.com, password=c2FyaXRoYXJhcGhlbA==, memberType=PRIMARY, role=USER, memberDOB=2025-05-16, mobileNumber=7736544685, address1=null, place=null, taluk=Kottayam, district=Kottayam, state=Kerala, pincode=688537)"
$ws.Range("G21").Value = "MemberModel(no=7, memberName=മനേഷ് മൈക്കിൽ, emailId=manesh.micheal@gmail.com, password=bWFuZXNoLm1pY2hlYWw=, memberType=PRIMARY, role=USER, memberDOB=2025-11-07, mobileNumber=9995227002, address1=null, place=null, taluk=Cherthala, district=Alappuzha, state=Kerala, pincode=688537)"
$ws.Range("G22").Value = "MemberModel(no=10, memberName=സരിത റഫയൽ, emailId=saritharaphel@gmail.com, password=c2FyaXRoYXJhcGhlbA==, memberType=PRIMARY, role=USER, memberDOB=2025-05-16, mobileNumber=7736544685, address1=null, place=null, taluk=Kottayam, district=Kottayam, state=Kerala, pincode=688537)"
